$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the values in D1 and E1, but keep their existing cell formatting/style
$ws.Range("D1:E1").ClearContents()
